# Applies updated values to the RLIe_GET_L3 worksheet (Sheet1) as described
# in the commit "added narratives to the findings". Only numeric cell
# values in columns F (lower) and G (upper) changed for several rows,
# plus one value in column G for row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("G5").Value = 0.78

$ws.Range("F18").Value = 0.5833333333333333
$ws.Range("G18").Value = 0.7060606060606061

$ws.Range("G19").Value = 0.6909090909090909

$ws.Range("G20").Value = 0.6909090909090909

$ws.Range("F21").Value = 0.5696969696969697
$ws.Range("G21").Value = 0.6909090909090909

$ws.Range("F22").Value = 0.8842105263157894

$ws.Range("F23").Value = 0.8789473684210526

$ws.Range("F24").Value = 0.8631578947368421

$ws.Range("F25").Value = 0.8631578947368421

$ws.Range("F27").Value = 0.7462686567164178

$ws.Range("G28").Value = 0.8746268656716418

$ws.Range("F29").Value = 0.7373134328358208
$ws.Range("G29").Value = 0.8746268656716418

$ws.Range("F31").Value = 0.76
$ws.Range("G31").Value = 0.888

$ws.Range("F32").Value = 0.7573333333333333
$ws.Range("G32").Value = 0.8853333333333333

$ws.Range("F33").Value = 0.7493333333333334
$ws.Range("G33").Value = 0.8773333333333333

$ws.Range("F38").Value = 0.8603174603174604

$ws.Range("F40").Value = 0.8444444444444444

$ws.Range("F41").Value = 0.8412698412698413
